# Crackpots (Atari 2600) - Plan.xlsx edit
# Updates the "multiplier" and "combo" achievement descriptions so they no
# longer use a parenthetical aside, per the new writing policy:
#   "Get a Nx multiplier! (...)"   ->  "Get a Nx multiplier! ... ."
#   "Squash N bugs in a combo! (...)" -> "Squash N bugs in a combo! ... ."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Achievements")

$ws.Range("F6").Value  = "Get a 2x multiplier! Beat all four bug waves to increase a mutliplier, losing a floor will decrease a mutliplier."
$ws.Range("F7").Value  = "Get a 3x multiplier! Beat all four bug waves to increase a mutliplier, losing a floor will decrease a mutliplier."
$ws.Range("F8").Value  = "Get a 4x multiplier! Beat all four bug waves to increase a mutliplier, losing a floor will decrease a mutliplier."
$ws.Range("F9").Value  = "Get a 5x multiplier! Beat all four bug waves to increase a mutliplier, losing a floor will decrease a mutliplier."

$ws.Range("F19").Value = "Squash 2 bugs in a combo! Combos end when all the pots have been reset."
$ws.Range("F20").Value = "Squash 3 bugs in a combo! Combos end when all the pots have been reset."
$ws.Range("F21").Value = "Squash 5 bugs in a combo! Combos end when all the pots have been reset."
$ws.Range("F22").Value = "Squash 7 bugs in a combo! Combos end when all the pots have been reset."

# Recalculate so dependent sheets (Checklist, Text) refresh their cached
# formula values that quote these achievement descriptions.
$excel.CalculateFullRebuild()

# Restore the sheet's remembered selection to match the saved state.
$ws.Activate()
$ws.Range("F23").Select()
